# "working on the complementarity"
# Merge the two smallest/duplicate marginal-cost blocks on each sheet
# (row 6 + row 7) into a single row 6: Pmax [MW] is summed, Bid price
# [$/MWh] is averaged. Row 7 is then removed from each sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Generation_investor")
$ws2 = $wb.Worksheets.Item("Generation_rival")

# --- Generation_rival: edit first, then drop row 7 ---
$ws2.Activate()
$ws2.Range("C6").Formula = "=350+310"
$ws2.Range("D6").Formula = "=(10.52+10.89)/2"
$ws2.Rows("7:7").Select()
$ws2.Rows("7:7").Delete()

# --- Generation_investor: edit, drop row 7, leave it the active sheet ---
$ws1.Activate()
$ws1.Range("C6").Formula = "=60+155"
$ws1.Range("D6").Formula = "=(26.11+10.52)/2"
$ws1.Rows("7:7").Delete()
$ws1.Range("H11").Select()
